# Add I0 and IF columns to the sheet
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Headers in I1 and J1, copying the style used by other header cells (e.g. H1)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Data values for rows 2-23
$data = @{
    2  = @(6, 7)
    3  = @(7, 7)
    4  = @(8, 8)
    5  = @(8, 8)
    6  = @(5, 6)
    7  = @(7, 7)
    8  = @(8, 8)
    9  = @(8, 8)
    10 = @(10, 10)
    11 = @(9, 9)
    12 = @(7, 7)
    13 = @(8, 8)
    14 = @(8, 8)
    15 = @(5, 5)
    16 = @(8, 8)
    17 = @(3, 3)
    18 = @(6, 6)
    19 = @(7, 7)
    20 = @(7, 7)
    21 = @(9, 9)
    22 = @(6, 6)
    23 = @(5, 5)
}

foreach ($row in $data.Keys) {
    $vals = $data[$row]
    $ws.Cells.Item($row, 9).Value = $vals[0]
    $ws.Cells.Item($row, 10).Value = $vals[1]
}
